$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13:131 down to 14:132.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new data record.
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = 'Macroferia Regional de Talca'
$ws.Range("C13").Value = 'Maule'
$ws.Range("D13").Value = '2022-07-08'
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 100112001
$ws.Range("G13").Value = 'Berenjena'
$ws.Range("H13").Value = 'Sin especificar'
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 8000
$ws.Range("L13").Value = 8000
$ws.Range("M13").Value = 8000
$ws.Range("N13").Value = '$/caja 50 unidades'
$ws.Range("O13").Value = 'Región de Arica y Parinacota'
$ws.Range("P13").Value = 160
$ws.Range("Q13").Value = 50
$ws.Range("R13").Value = 'Hortaliza'
